$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.157014012336731
$ws.Range("B1").Value = 1.380603790283203
$ws.Range("C1").Value = 1.141376256942749
$ws.Range("D1").Value = 1.129527807235718
$ws.Range("E1").Value = 1.161216974258423
